# (C)    updated/corrected a few lifts, cleaned up doc
#
# - Swap the "Lower (9)" (Day 2) and "Arms (9)" (Day 3) blocks: the Arms
#   exercises now come first (Day 2), followed by the Lower exercises (Day 3).
# - Fix a couple of lift names:
#     * "(Cable Fly Seated)"            -> "(Selectorized Fly)"
#     * "(Barbell Raise)"               -> "(Barbbell Raise)"
#     * "(Dumbbell Biceps Curl (Underhand Grip)"  -> "...Grip))" (closing paren)
#     * "(Machine Revese Fky)"          -> "(Machine Revese Fly)"
# - Bump the revision note in I45 from "(11/23/16) r5" to "(11/24/16) r6".
# - Leave the selection on D19.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- small text corrections that are not part of the block swap ---
$ws.Range("D7").Value2  = "(Selectorized Row)"
$ws.Range("D8").Value2  = "(Barbbell Raise)"
$ws.Range("D12").Value2 = "(Selectorized Fly)"

# --- capture the current ("Lower (9)") block, rows 17-26, columns B & D ---
$lowerB = @{}
$lowerD = @{}
for ($r = 17; $r -le 26; $r++) {
    $lowerB[$r] = $ws.Cells.Item($r, 2).Value2
    $lowerD[$r] = $ws.Cells.Item($r, 4).Value2
}

# --- capture the current ("Arms (9)") block, rows 28-37, columns B & D ---
$armsB = @{}
$armsD = @{}
for ($r = 28; $r -le 37; $r++) {
    $armsB[$r] = $ws.Cells.Item($r, 2).Value2
    $armsD[$r] = $ws.Cells.Item($r, 4).Value2
}

# fix two typos while they are in hand, before they get written back out
$armsD[29] = "(Dumbbell Biceps Curl (Underhand Grip))"
$armsD[35] = "(Machine Revese Fly)"

# --- write the Arms block into rows 17-26 ---
for ($i = 0; $i -le 9; $i++) {
    $destRow = 17 + $i
    $srcRow  = 28 + $i
    $ws.Cells.Item($destRow, 2).Value2 = $armsB[$srcRow]
    $ws.Cells.Item($destRow, 4).Value2 = $armsD[$srcRow]
}

# --- write the Lower block into rows 28-37 ---
for ($i = 0; $i -le 9; $i++) {
    $destRow = 28 + $i
    $srcRow  = 17 + $i
    $ws.Cells.Item($destRow, 2).Value2 = $lowerB[$srcRow]
    $ws.Cells.Item($destRow, 4).Value2 = $lowerD[$srcRow]
}

# --- bump the revision marker ---
$ws.Range("I45").Value2 = "(11/24/16) r6"

# --- leave the selection on D19, matching the saved file ---
$ws.Activate()
$ws.Range("D19").Select()
